$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.128.17"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "'1.666.86"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "'209.53"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").Value = "'0.5205"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "'0.2616"
$ws.Range("E8").Value = "  -3.43%  "
$ws.Range("D9").Value = "'0.06316"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "'21.11"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").Value = "'0.07535"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "'1.679.22"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'4.424"
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("E14").Value = "  -5.22%  "
$ws.Range("D15").Value = "'66.25"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'0.000007929"
$ws.Range("E16").Value = "  -5.16%  "
$ws.Range("D17").Value = "'26.142.55"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'4.722"
$ws.Range("D20").Value = "'186.35"
$ws.Range("E20").Value = "  -2.42%  "
$ws.Range("E21").Value = "  -5.29%  "
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "'149.00"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'0.1247"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "'7.476"
$ws.Range("E26").Value = "  -4.41%  "
$ws.Range("D27").Value = "'15.82"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'0.06370"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'1.351"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "'1.273"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").Value = "'3.494"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").Value = "'3.406"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "'1.002"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("D36").Value = "'0.6005"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "'1.109.46"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").Value = "'6.094"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").Value = "'0.01612"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'0.8649"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "'100.08"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'1.819.63"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "'0.9974"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "'8.024"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'0.4247"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("D51").Value = "'5.913"
$ws.Range("E51").Value = "  -2.18%  "
